$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 0.9823068853471902
$ws.Range("D4").Value = 0.2254840312804062
$ws.Range("E4").Value = 1.197280435915919
$ws.Range("F4").Value = 0.2254840312804062
$ws.Range("G4").Value = 15.10956701022704
$ws.Range("H4").Value = 69.04212945962217
$ws.Range("I4").Value = 15.84830353015079
$ws.Range("J4").Value = 84.15169646984921
$ws.Range("K4").Value = 15.84830353015079

# Row 5
$ws.Range("D5").Value = 0.9839436444806937
$ws.Range("F5").Value = 0.9839436444806937
$ws.Range("G5").Value = 3.560091977063498
$ws.Range("H5").Value = 83.02406258426124
$ws.Range("I5").Value = 13.41584543867526
$ws.Range("J5").Value = 86.58415456132474
$ws.Range("K5").Value = 13.41584543867526

# Row 20 - values AND style change (s=5 -> s=7, i.e. #,##0.000 format)
$ws.Range("C20").Value = 0.04298276364640884
$ws.Range("D20").Value = 0.03504234798678382
$ws.Range("E20").Value = 0.04298276364640884
$ws.Range("F20").Value = 0.03504234798678382
$ws.Range("C20:F20").NumberFormat = "#,##0.000"
$ws.Range("H20").Value = 55.08837186735091
$ws.Range("I20").Value = 44.91162813264909
$ws.Range("J20").Value = 55.08837186735091
$ws.Range("K20").Value = 44.91162813264909

# Row 22
$ws.Range("B22").Value = 0.02678414
$ws.Range("C22").Value = 0.02084913
$ws.Range("D22").Value = 0.00601671
$ws.Range("E22").Value = 0.04763327000000001
$ws.Range("F22").Value = 0.00601671
$ws.Range("G22").Value = 49.92385831271513
$ws.Range("H22").Value = 38.86139379735091
$ws.Range("I22").Value = 11.21474788993398
$ws.Range("J22").Value = 88.78525211006604
$ws.Range("K22").Value = 11.21474788993398

# Row 24
$ws.Range("B24").Value = 18.50865755204852
$ws.Range("C24").Value = 33.60113412567404
$ws.Range("D24").Value = 17.34828565888902
$ws.Range("E24").Value = 52.10979167772258
$ws.Range("F24").Value = 17.34828565888902
$ws.Range("G24").Value = 26.64723565892969
$ws.Range("H24").Value = 48.37613624522655
$ws.Range("I24").Value = 24.97662809584376
$ws.Range("J24").Value = 75.02337190415625
$ws.Range("K24").Value = 24.97662809584376
